# Daily price-sheet refresh: a new row is inserted just under the header
# for the latest date, pushing every existing date row down by one. The
# price figures (col B/C/D) stay the same as the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Insert()

# Leading apostrophe forces plain text entry so "2026-01-27" is not
# auto-converted into a date serial (matching the rest of column A, which
# stores the dates as literal text). Re-copy the sibling cell's style so
# the text-entry doesn't leave a stray "quote prefix" formatting flag on
# this cell.
$ws.Range("A2").Value = "'2026-01-27"
$ws.Range("A2").Style = $ws.Range("A3").Style

$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
